$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("pop", "births", "deaths", "__axes__")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("A1").Value = "country"
}
